# Insert a new data row before row 498 (shifting existing rows 498-604 down
# to 499-605) and populate the new row with a new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 498; this pushes the old row 498 (and
# everything below it) down to row 499, and all formatting/values carry
# along automatically.
$ws.Rows.Item(498).Insert()

# Populate the constant columns (identical for every record in this block)
$ws.Cells.Item(498, 1).Value = 3                               # A: Mercado ID
$ws.Cells.Item(498, 2).Value = "Femacal de La Calera"          # B: Mercado
$ws.Cells.Item(498, 3).Value = "Coquimbo"                      # C: Region
$ws.Cells.Item(498, 4).Value = 45173                           # D: Fecha
$ws.Cells.Item(498, 5).Value = 5                                # E: Codreg
$ws.Cells.Item(498, 6).Value = 100112012                        # F: Categoria ID
$ws.Cells.Item(498, 7).Value = "Espinaca"                       # G: Categoria
$ws.Cells.Item(498, 8).Value = "Sin especificar"                # H: Variedad
$ws.Cells.Item(498, 9).Value = "Primera"                        # I: Calidad
$ws.Cells.Item(498, 10).Value = 110                             # J: Volumen
$ws.Cells.Item(498, 11).Value = 4000                            # K: Precio minimo
$ws.Cells.Item(498, 12).Value = 4000                            # L: Precio maximo
$ws.Cells.Item(498, 13).Value = 4000                            # M: Precio promedio ponderado
$ws.Cells.Item(498, 14).Value = "`$/docena de atados (3 kilos)" # N: Unidad de comercializacion
$ws.Cells.Item(498, 15).Value = "Provincia de Quillota"         # O: Origen
$ws.Cells.Item(498, 16).Value = 1333                            # P: Precio $/Kg
$ws.Cells.Item(498, 17).Value = 3                               # Q: Kg o Unidades
$ws.Cells.Item(498, 18).Value = "Hortaliza"                     # R: Clasificacion

# Make sure the date cell keeps the date/time number format used by the
# rest of column D.
$ws.Cells.Item(498, 4).NumberFormat = $ws.Cells.Item(499, 4).NumberFormat
